# "add pass test mark and add ranking function"
#
# For every per-college data sheet (index 2..8 — the seven "1.x.x.x ..." /
# "2.x.x.x ..." / "3.x.x.x ..." indicator sheets that sit between the
# "小結" summary sheet and the end of the tab strip):
#   - the "999" pass/fail test mark that used to sit on the 校均值/校加總
#     row (G2) is removed
#   - a "999" mark is written onto the last row instead (G13), which is
#     where the new ranking-column helper value belongs
#   - the remembered cell selection is moved from F7 to D14
#
# The "小結" (summary) sheet additionally gets an explicit print/page setup
# (A4, portrait) recorded on it.

$wb = $excel.ActiveWorkbook

# --- 小結 (summary) sheet: record explicit page setup -----------------
$summary = $wb.Worksheets.Item(1)
$summary.PageSetup.PaperSize = 9        # xlPaperA4
$summary.PageSetup.Orientation = 1      # xlPortrait

# --- the seven indicator sheets ---------------------------------------
for ($i = 2; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Move the "999" pass-test mark from the average row (G2) down to the
    # new ranking row (G13).
    $ws.Range("G2").ClearContents() | Out-Null
    $ws.Range("G13").Value = 999

    # Remember the new selected cell used while reviewing the ranking
    # column.
    $ws.Range("D14").Select() | Out-Null
}

Write-Output "done"
